$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Análisis")

# Update values in column D (T(Práctico) (ns) data used by chart1)
$ws.Range("D4").Value = 122500
$ws.Range("D5").Value = 1741300
$ws.Range("D6").Value = 29079800

# Update values in column I (used by chart4)
$ws.Range("I4").Value = 1145900
$ws.Range("I5").Value = 3296900
$ws.Range("I6").Value = 29168300

# Update the sheet view: reset scrolled topLeftCell and change the selection
$ws.Activate()
$ws.Range("E5").Select()
